# 20-Feb-2021 end of day update for the "Buku KAS HARIAN" petty cash book.
# Target worksheet is the first tab (internally sheet1.xml), which holds
# the daily cash-book ledger with columns: A=Tgl, B=Keterangan, C=Debit,
# D=Credit, E=Saldo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 38 - Wages Expense day: add a Credit (D) entry.
$ws.Range("D38").Formula = "=60000+300000"

# Row 39 - A/R: add a Debit (C) entry.
$ws.Range("C39").Formula = "=145000+19813000"

# Row 40 - TRANSFER BCA: add a Credit (D) entry.
$ws.Range("D40").Formula = "=145000+1940000+120000"

# Row 41 - new line item: SALES - cash/retail, with a Debit (C) entry.
$ws.Range("B41").Value = "SALES - cash/retail"
$ws.Range("C41").Formula = "=934975+28543525-19813000"

# Row 42 - new line item: SELISIH - kurang, with a plain Credit (D) value.
$ws.Range("B42").Value = "SELISIH - kurang"
$ws.Range("D42").Value = 20000

# Row 43 - new line item: SETOR KE BANK, with a plain Credit (D) value.
$ws.Range("B43").Value = "SETOR KE BANK"
$ws.Range("D43").Value = 28000000

# Row 44 - new day entry: 20-Feb-2021, Wages Expense.
$ws.Range("A44").Value = 44247
$ws.Range("A44").NumberFormat = $ws.Range("A38").NumberFormat
$ws.Range("B44").Value = "Wages Expense"

# Update the active selection to reflect the end-of-day view position
# (scrolled down towards the newly entered rows, frozen header still
# showing rows 1-2).
$ws.Activate()
$ws.Range("B45").Select()
